$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure each updated cell keeps a text format so the numeric-looking
# strings (e.g. "281.31") and percentages (e.g. "6.22%") are stored as
# text, matching the original inlineStr/text content of the workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "281.31"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.22%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.06%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.934"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.98%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06393"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.91%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.978"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.62%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.348"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8838"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9545"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.12%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1496"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.52%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05197"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.40%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07463"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.15%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03115"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.86%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09028"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.01%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001584"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.70%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006324"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.07%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005987"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.03%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.506"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.65%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.298"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.69%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3094"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.75%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1289"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.60%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.930"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.75%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04333"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.72%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001174"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.31%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003667"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.57%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001196"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.44%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001690"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04102"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.46%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006643"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "58.23%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1175"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.48%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002351"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "11.62%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.32%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005244"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.64%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.08%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "811.99%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02245"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-8.30%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.08%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.15%"
